$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")  # "Hoja1" == $wb.ActiveSheet (tabSelected sheet)

# Row 1
$ws.Range("A1").Value = "FALTA"

# Row 3 - headers (bold)
$ws.Range("A3").Value = "DOCUMENTO"
$ws.Range("B3").Value = "DESCRIPCION"
$ws.Range("A3:B3").Font.Bold = $true

# Row 4
$ws.Range("A4").Value = "Documento de requisitos"
$ws.Range("B4").Value = "pantallazos de toda la parte servidora"

# Row 5
$ws.Range("A5").Value = "Documento de plan de proyecto"
$ws.Range("B5").Value = "IP servidor de desarrollo"

# Row 6
$ws.Range("A6").Value = "analisis del sistema"
$ws.Range("B6").Value = "Explicacion de las tablas que faltan"

# Column widths (best-fit for "DOCUMENTO"/"DESCRIPCION" columns, ~29.71 / ~34.57 chars)
$ws.Columns.Item(1).ColumnWidth = 28.82
$ws.Columns.Item(2).ColumnWidth = 33.65

# Selection
$ws.Range("A9").Select() | Out-Null
